# Insert a new data row at row 942 (pushing existing rows 942-996 down to 943-997)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("942").Insert()

$ws.Range("A942").Value = 5
$ws.Range("B942").Value = "Macroferia Regional de Talca"
$ws.Range("C942").Value = "Maule"
$ws.Range("D942").Value = 45267
$ws.Range("E942").Value = 7
$ws.Range("F942").Value = 100114001
$ws.Range("G942").Value = "Papa"
$ws.Range("H942").Value = "Asterix"
$ws.Range("I942").Value = "1a nueva(o)"
$ws.Range("J942").Value = 2300
$ws.Range("K942").Value = 16000
$ws.Range("L942").Value = 18000
$ws.Range("M942").Value = 17304
$ws.Range("N942").Value = "`$/saco 25 kilos"
$ws.Range("O942").Value = "Región del Maule"
$ws.Range("P942").Value = 692
$ws.Range("Q942").Value = 25
$ws.Range("R942").Value = "Hortaliza"
